$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing "climate_temperate" row (currently row 19)
# to hold the new "climate_title" / "Climate" key-value pair.
$ws.Rows.Item(19).Insert()
$ws.Cells.Item(19, 1).Value = "climate_title"
$ws.Cells.Item(19, 2).Value = "Climate"
$ws.Cells.Item(19, 2).WrapText = $true

# Insert a new row before the existing "region_NA" row (now shifted to row 21)
# to hold the new "region_title" / "Region" key-value pair.
$ws.Rows.Item(21).Insert()
$ws.Cells.Item(21, 1).Value = "region_title"
$ws.Cells.Item(21, 2).Value = "Region"
$ws.Cells.Item(21, 2).WrapText = $true

$ws.Range("A21").Select()
